# Updates the unrealized gain/loss worksheet with refreshed market-price figures
# (investment cost, market value, P&L, P&L%, trial price, trial market value)
# for each holding row, plus the recomputed "小計" (subtotal) row 21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell="D2"; Value=197},
    @{Cell="E2"; Value=44.9},
    @{Cell="F2"; Value=8851},
    @{Cell="G2"; Value=9978},
    @{Cell="H2"; Value=1127},
    @{Cell="I2"; Value=0.127},
    @{Cell="K2"; Value=10047},
    @{Cell="E3"; Value=49.64},
    @{Cell="F3"; Value=2881},
    @{Cell="G3"; Value=2640},
    @{Cell="H3"; Value=-241},
    @{Cell="I3"; Value=-0.084},
    @{Cell="J3"; Value=46},
    @{Cell="K3"; Value=2668},
    @{Cell="G4"; Value=4550},
    @{Cell="H4"; Value=-571},
    @{Cell="I4"; Value=-0.112},
    @{Cell="J4"; Value=139.5},
    @{Cell="K4"; Value=4604},
    @{Cell="G5"; Value=2508},
    @{Cell="H5"; Value=-375},
    @{Cell="I5"; Value=-0.13},
    @{Cell="J5"; Value=81.8},
    @{Cell="K5"; Value=2536},
    @{Cell="G6"; Value=137942},
    @{Cell="H6"; Value=19207},
    @{Cell="I6"; Value=0.162},
    @{Cell="J6"; Value=1155},
    @{Cell="K6"; Value=138600},
    @{Cell="G7"; Value=4008},
    @{Cell="H7"; Value=822},
    @{Cell="I7"; Value=0.258},
    @{Cell="J7"; Value=101},
    @{Cell="K7"; Value=4040},
    @{Cell="G8"; Value=242767},
    @{Cell="H8"; Value=44117},
    @{Cell="I8"; Value=0.222},
    @{Cell="J8"; Value=66.8},
    @{Cell="K8"; Value=243954},
    @{Cell="G9"; Value=51975},
    @{Cell="H9"; Value=-8518},
    @{Cell="I9"; Value=-0.141},
    @{Cell="J9"; Value=79.2},
    @{Cell="K9"; Value=52351},
    @{Cell="G10"; Value=2676},
    @{Cell="H10"; Value=-283},
    @{Cell="I10"; Value=-0.096},
    @{Cell="J10"; Value=41.6},
    @{Cell="K10"; Value=2704},
    @{Cell="G11"; Value=39837},
    @{Cell="H11"; Value=-2189},
    @{Cell="I11"; Value=-0.052},
    @{Cell="J11"; Value=82.5},
    @{Cell="K11"; Value=40012},
    @{Cell="G12"; Value=29718},
    @{Cell="H12"; Value=712},
    @{Cell="I12"; Value=0.025},
    @{Cell="J12"; Value=31.6},
    @{Cell="K12"; Value=29862},
    @{Cell="G13"; Value=5850},
    @{Cell="H13"; Value=-217},
    @{Cell="I13"; Value=-0.036},
    @{Cell="J13"; Value=84.1},
    @{Cell="K13"; Value=5887},
    @{Cell="G14"; Value=6702},
    @{Cell="H14"; Value=1374},
    @{Cell="I14"; Value=0.258},
    @{Cell="J14"; Value=36.95},
    @{Cell="K14"; Value=6762},
    @{Cell="G15"; Value=2191},
    @{Cell="H15"; Value=-842},
    @{Cell="I15"; Value=-0.278},
    @{Cell="J15"; Value=48.2},
    @{Cell="K15"; Value=2217},
    @{Cell="G16"; Value=5542},
    @{Cell="H16"; Value=-420},
    @{Cell="I16"; Value=-0.07},
    @{Cell="J16"; Value=155.5},
    @{Cell="K16"; Value=5598},
    @{Cell="G17"; Value=7346},
    @{Cell="H17"; Value=671},
    @{Cell="I17"; Value=0.101},
    @{Cell="J17"; Value=823},
    @{Cell="K17"; Value=7407},
    @{Cell="D18"; Value=22},
    @{Cell="E18"; Value=144.73},
    @{Cell="F18"; Value=3187},
    @{Cell="G18"; Value=3087},
    @{Cell="H18"; Value=-100},
    @{Cell="I18"; Value=-0.031},
    @{Cell="J18"; Value=142.5},
    @{Cell="K18"; Value=3135},
    @{Cell="E19"; Value=212.92},
    @{Cell="F19"; Value=2770},
    @{Cell="G19"; Value=2482},
    @{Cell="H19"; Value=-288},
    @{Cell="I19"; Value=-0.104},
    @{Cell="J19"; Value=193},
    @{Cell="K19"; Value=2509},
    @{Cell="G20"; Value=3008},
    @{Cell="H20"; Value=-174},
    @{Cell="I20"; Value=-0.055},
    @{Cell="J20"; Value=121.5},
    @{Cell="K20"; Value=3038},
    @{Cell="F21"; Value=510995},
    @{Cell="G21"; Value=564807},
    @{Cell="H21"; Value=53812},
    @{Cell="I21"; Value=0.105},
    @{Cell="K21"; Value=567931}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
